$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 201, shifting existing rows 201-246 down to 202-247
$ws.Range("A201").EntireRow.Insert()

# Populate the new row 201 with the new weekly record
$ws.Range("A201").Value = 4
$ws.Range("B201").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C201").Value = "Los Lagos"
$ws.Range("D201").Value = 44722
$ws.Range("E201").Value = 10
$ws.Range("F201").Value = 100112032
$ws.Range("G201").Value = "Zapallo italiano"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 200
$ws.Range("K201").Value = 15000
$ws.Range("L201").Value = 15000
$ws.Range("M201").Value = 15000
$ws.Range("N201").Value = '$/caja 50 unidades'
$ws.Range("O201").Value = "Región de Arica y Parinacota"
$ws.Range("P201").Value = 300
$ws.Range("Q201").Value = 50
$ws.Range("R201").Value = "Hortaliza"
